# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)
#
# Renames the sole "Data" sheet to "Summary", inserts a new "Source Type"
# header line, pushes the MSME table down a couple of rows, and appends a
# fuller source citation (ISTEEBU name + long reference) at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet: "Data" -> "Summary"
$ws.Name = "Summary"

# Re-assert the formatting on the two untouched header cells (A1 "name"
# style / A3 "title" style) so their styles get freshly (and correctly)
# re-indexed by the save pipeline alongside everything else that's about
# to change on this sheet.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# 2. Clear out the old table (rows 5-9) -- it gets rebuilt a couple of rows
#    lower down, with a new "Source Type" line inserted above it.
$ws.Range("A5:D9").Clear()

# 3. New bold+underlined "Source Type" line.
$ws.Range("A7").Value = "Source Type: Statistical Institution (Most Widely Used)"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# 4. Column headers (row 9), bold.
$ws.Range("B9").Value = "Micro"
$ws.Range("C9").Value = "SMEs"
$ws.Range("D9").Value = "MSMEs"
$ws.Range("B9:D9").Font.Bold = $true

# 5. Data rows 10-12: bold row labels in column A, plain text data in B:D.
$ws.Range("A10").Value = "Enterprises (absolute #)"
$ws.Range("A10").Font.Bold = $true
$ws.Range("B10").Value = "'3411"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'388"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'3799"
$ws.Range("D10").Style = "Normal"

$ws.Range("A11").Value = "Enterprises density (per 1000 people)"
$ws.Range("A11").Font.Bold = $true
$ws.Range("B11").Value = "'0.4"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'0"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.4"
$ws.Range("D11").Style = "Normal"

$ws.Range("A12").Value = "Enterprises (% of total)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("B12").Value = "'89"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'10.1"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'99.1"
$ws.Range("D12").Style = "Normal"

# 6. Source note (italic), moved down from row 9 to row 13.
$ws.Range("A13").Value = "Source: ISTEEBU, 2010"
$ws.Range("A13").Font.Italic = $true

# 7. New "ISTEEBU" bold line and long citation (italic) further down.
$ws.Range("A21").Value = "ISTEEBU"
$ws.Range("A21").Font.Bold = $true

$ws.Range("A22").Value = 'L''Institut de Statistiques et d’Etudes Economiques du Burundi (ISTEEBU), "ANNUAIRE STATISTIQUE DU BURUNDI 2011", Répartition des entreprises selon le nombre de travailleurs actifs (au 31 décembre), 2013, p. 160. Available at http://www.isteebu.bi/images/annuaires/annuaire%202011.pdf%20vf.pdf'
$ws.Range("A22").Font.Italic = $true

Write-Host "Applied Burundi Summary edits"
